# Updated documentation. Added new Integration test written using Protractor.
#
# The "Test cases" sheet grows from 9 to 17 test-case rows (sheet rows 2-18):
#   - Row 2/3 (Test Case 1/2) gain "Actual Results"/"Pass-Fail" entries (cols F/G)
#     that were previously blank, and row 3's "Expected results" text is replaced.
#   - Row 4 (Test Case 3, previously row 3 in the old sheet) gets new F/G results.
#   - Test Case 4-10 (previously the "sort" rows 4-10) are rewritten with the new
#     Protractor results and two brand new cases (Test Case 9 "Race and Comfort"
#     and Test Case 10 "Endurance Race and Comfort") are appended.
#   - The old filter rows (previously Test Case 4-9, "Endurance"/"Race"/.../
#     "Race, Comfort") move down to Test Case 11-17 (sheet rows 12-18), each
#     gaining F/G result columns, plus one brand new case (Test Case 17,
#     "Endurance, Race, Comfort") using the original 3-step refresh-page flow.
#   - The two leftover placeholder strings ("???" and the "I don't see any sort
#     functionality" note) are dropped along with the row that held them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------------
# 1) Grow the formatted table down to row 18.
#    Row 10 already has the correct borders / wrap-text / font styling
#    (style "2" for columns A,C,D,E,F,G and style "3" for column B); clone
#    it across the 8 new rows so they pick up the same look instead of the
#    COM default style.
# ------------------------------------------------------------------------
$ws.Range("A10:G10").Copy()
$ws.Range("A11:G18").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 4's "B" cell (the sort-related scenario) keeps the un-bordered style
# that the original sheet used for that particular scenario text (same
# style as it had one row up, before the new rows were inserted above it).
$ws.Range("B3").Copy()
$ws.Range("B4").PasteSpecial(-4122)
$excel.CutCopyMode = $false


# --- Row 2 (Test Case 1) ------------------------------
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "I would like to see a list/grid of bikes based on /app/bikes.json"
$ws.Range("C2").Value = "Open app"
$ws.Range("D2").Value = "/app/bikes.json"
$ws.Range("E2").Value = @"
All bikes from /app/bikes.json
should be visible on screen
(should be present in the HTML document).
"@
$ws.Range("F2").Value = "7 bikes detected from JSON file. 7 bikes detected from app controller."
$ws.Range("G2").Value = "Pass"
$ws.Rows.Item(2).RowHeight = 105

# --- Row 3 (Test Case 2) ------------------------------
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "I would like to see a list/grid of bikes based on /app/bikes.json"
$ws.Range("C3").Value = "Open app"
$ws.Range("D3").Value = "/app/bikes.json"
$ws.Range("E3").Value = @"
The layout of the screen should be defined by a layout that presents as a grid, but resizes to a list when the screen narrow.
The following class should be expected on the element that generates the bikelist:
<div ng-repeat="p in filteredProducts" class="col-sm-4">
"@
$ws.Range("F3").Value = "class name confirmed as col-sm-4"
$ws.Range("G3").Value = "Pass"
$ws.Rows.Item(3).RowHeight = 150

# --- Row 4 (Test Case 3) ------------------------------
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "I would like to see a the name, an image, a description and the class for each bike"
$ws.Range("C4").Value = "Open app"
$ws.Range("D4").Value = "/app/bikes.json"
$ws.Range("E4").Value = @"
For each bike, the following elements must be present:
Name, Image, Description, Class 
"@
$ws.Range("F4").Value = "The definition of all elements has  been confirmed"
$ws.Range("G4").Value = "Pass"
$ws.Rows.Item(4).RowHeight = 75

# --- Row 5 (Test Case 4) ------------------------------
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "I would like to be able to sort the bikes into a custom order based on class"
$ws.Range("C5").Value = @"
1.Open app
2.Filter By Endurance  
3.Get products on display
4.Toggle filter off
5. Recount products

"@
$ws.Range("D5").Value = "/app/bikes.json"
$ws.Range("E5").Value = "6 bikes should show"
$ws.Range("F5").Value = "Correct"
$ws.Range("G5").Value = "Pass"
$ws.Rows.Item(5).RowHeight = 90

# --- Row 6 (Test Case 5) ------------------------------
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "I would like to be able to sort the bikes into a custom order based on class"
$ws.Range("C6").Value = @"
1.Open app
2.Filter By Race
3.Get products on display
4.Toggle filter off
5. Recount products
"@
$ws.Range("D6").Value = "/app/bikes.json"
$ws.Range("E6").Value = "4 bikes should show"
$ws.Range("F6").Value = "Correct"
$ws.Range("G6").Value = "Pass"
$ws.Rows.Item(6).RowHeight = 75

# --- Row 7 (Test Case 6) ------------------------------
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = "I would like to be able to sort the bikes into a custom order based on class"
$ws.Range("C7").Value = @"
1.Open app
2.Filter By Comfort
3.Get products on display
4.Toggle filter off
5. Recount products
"@
$ws.Range("D7").Value = "/app/bikes.json"
$ws.Range("E7").Value = "4 bikes should show"
$ws.Range("F7").Value = "Correct"
$ws.Range("G7").Value = "Pass"
$ws.Rows.Item(7).RowHeight = 75

# --- Row 8 (Test Case 7) ------------------------------
$ws.Range("A8").Value = 7
$ws.Range("B8").Value = "I would like to be able to sort the bikes into a custom order based on class"
$ws.Range("C8").Value = @"
1.Open app
2.Filter By Enduance and Comfort
3.Get products on display
4.Toggle filter off
5. Recount products
"@
$ws.Range("D8").Value = "/app/bikes.json"
$ws.Range("E8").Value = "3 bikes should show"
$ws.Range("F8").Value = "Correct"
$ws.Range("G8").Value = "Pass"
$ws.Rows.Item(8).RowHeight = 75

# --- Row 9 (Test Case 8) ------------------------------
$ws.Range("A9").Value = 8
$ws.Range("B9").Value = "I would like to be able to sort the bikes into a custom order based on class"
$ws.Range("C9").Value = @"
1.Open app
2.Filter By Endurance and Race
3.Get products on display
4.Toggle filter off
5. Recount products
"@
$ws.Range("D9").Value = "/app/bikes.json"
$ws.Range("E9").Value = "4 bikes should show"
$ws.Range("F9").Value = "Correct"
$ws.Range("G9").Value = "Pass"
$ws.Rows.Item(9).RowHeight = 75

# --- Row 10 (Test Case 9) ------------------------------
$ws.Range("A10").Value = 9
$ws.Range("B10").Value = "I would like to be able to sort the bikes into a custom order based on class"
$ws.Range("C10").Value = @"
1.Open app
2.Filter By Race and Comfort
3.Get products on display
4.Toggle filter off
5. Recount products
"@
$ws.Range("D10").Value = "/app/bikes.json"
$ws.Range("E10").Value = "1 bikes should show"
$ws.Range("F10").Value = "Correct"
$ws.Range("G10").Value = "Pass"
$ws.Rows.Item(10).RowHeight = 75

# --- Row 11 (Test Case 10) ------------------------------
$ws.Range("A11").Value = 10
$ws.Range("B11").Value = "I would like to be able to sort the bikes into a custom order based on class"
$ws.Range("C11").Value = @"
1.Open app
2.Filter By Endurance Race and Comfort
3.Get products on display
4.Toggle filter off
5. Recount products
"@
$ws.Range("D11").Value = "/app/bikes.json"
$ws.Range("E11").Value = "1 bikes should show"
$ws.Range("F11").Value = "Correct"
$ws.Range("G11").Value = "Pass"
$ws.Rows.Item(11).RowHeight = 90

# --- Row 12 (Test Case 11) ------------------------------
$ws.Range("A12").Value = 11
$ws.Range("B12").Value = "I would like my custom order to be saved and not change when I refresh the page"
$ws.Range("C12").Value = @"
1.Open app
2.Filter By Endurance  
3.Refresh page
e.Get products on display
"@
$ws.Range("D12").Value = "/app/bikes.json"
$ws.Range("E12").Value = "The bike list should only show bikes that whose class includes Endurance."
$ws.Range("F12").Value = "Correct"
$ws.Range("G12").Value = "Pass"
$ws.Rows.Item(12).RowHeight = 60

# --- Row 13 (Test Case 12) ------------------------------
$ws.Range("A13").Value = 12
$ws.Range("B13").Value = "I would like my custom order to be saved and not change when I refresh the page"
$ws.Range("C13").Value = @"
1.Open app
2.Filter By Race 
3.Refresh page
4.Get products on display
"@
$ws.Range("D13").Value = "/app/bikes.json"
$ws.Range("E13").Value = "The bike list should only show bikes that whose class includes Race."
$ws.Range("F13").Value = "Correct"
$ws.Range("G13").Value = "Pass"
$ws.Rows.Item(13).RowHeight = 60

# --- Row 14 (Test Case 13) ------------------------------
$ws.Range("A14").Value = 13
$ws.Range("B14").Value = "I would like my custom order to be saved and not change when I refresh the page"
$ws.Range("C14").Value = @"
1.Open app
2.Filter By Comfort
3.Refresh page
4.Get products on display
"@
$ws.Range("D14").Value = "/app/bikes.json"
$ws.Range("E14").Value = "The bike list should only show bikes that whose class includes Comfort."
$ws.Range("F14").Value = "Correct"
$ws.Range("G14").Value = "Pass"
$ws.Rows.Item(14).RowHeight = 60

# --- Row 15 (Test Case 14) ------------------------------
$ws.Range("A15").Value = 14
$ws.Range("B15").Value = "I would like my custom order to be saved and not change when I refresh the page"
$ws.Range("C15").Value = @"
1.Open app
2.Filter By Endurance, Race
3.Refresh page
4.Get products on display
"@
$ws.Range("D15").Value = "/app/bikes.json"
$ws.Range("E15").Value = "The bike list should only show bikes that whose class includes Endurance, Race."
$ws.Range("F15").Value = "Correct"
$ws.Range("G15").Value = "Pass"
$ws.Rows.Item(15).RowHeight = 60

# --- Row 16 (Test Case 15) ------------------------------
$ws.Range("A16").Value = 15
$ws.Range("B16").Value = "I would like my custom order to be saved and not change when I refresh the page"
$ws.Range("C16").Value = @"
1.Open app
2.Filter By Endurance, Comfort
3.Refresh page
4.Get products on display
"@
$ws.Range("D16").Value = "/app/bikes.json"
$ws.Range("E16").Value = "The bike list should only show bikes that whose class includes Endurance, Comfort."
$ws.Range("F16").Value = "Correct"
$ws.Range("G16").Value = "Pass"
$ws.Rows.Item(16).RowHeight = 60

# --- Row 17 (Test Case 16) ------------------------------
$ws.Range("A17").Value = 16
$ws.Range("B17").Value = "I would like my custom order to be saved and not change when I refresh the page"
$ws.Range("C17").Value = @"
1.Open app
2.Filter By Race, Comfort
3.Refresh page
4.Get products on display
"@
$ws.Range("D17").Value = "/app/bikes.json"
$ws.Range("E17").Value = "The bike list should only show bikes that whose class includes Race, Comfort."
$ws.Range("F17").Value = "Correct"
$ws.Range("G17").Value = "Pass"
$ws.Rows.Item(17).RowHeight = 60

# --- Row 18 (Test Case 17) ------------------------------
$ws.Range("A18").Value = 17
$ws.Range("B18").Value = "I would like my custom order to be saved and not change when I refresh the page"
$ws.Range("C18").Value = @"
1.Open app
2.Filter By Endurance, Race, Comfort
3.Refresh page
4.Get products on display
"@
$ws.Range("D18").Value = "/app/bikes.json"
$ws.Range("E18").Value = "The bike list should only show bikes that whose class includes Race, Comfort."
$ws.Range("F18").Value = "Correct"
$ws.Range("G18").Value = "Pass"
$ws.Rows.Item(18).RowHeight = 75

# ------------------------------------------------------------------------
# Final cursor position, matching the saved selection in the workbook.
# ------------------------------------------------------------------------
$ws.Range("B25").Select()
